$d = $word.ActiveDocument

$bullet = [char]0x2022   # "•"

# Locate the "KEY ACHIEVEMENTS AND IMPACT" section heading. The section
# layout is:
#   <heading>         KEY ACHIEVEMENTS AND IMPACT
#   <sub-heading>     Impact
#   <bullet 1..6>      the six achievement bullets
# We operate by paragraph index (relative to the heading) rather than by
# text match, because some of the bullet text is duplicated earlier in the
# document (e.g. under PROFESSIONAL EXPERIENCE).
$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*KEY ACHIEVEMENTS AND IMPACT*") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq -1) {
    throw "Could not find KEY ACHIEVEMENTS AND IMPACT heading"
}

$bulletBase = $headingIndex + 1   # "Impact" sub-heading paragraph index

# Rewrite the first three bullets in place, preserving paragraph/run
# formatting (only the text content changes).
$d.Paragraphs($bulletBase + 1).Range.Text = $bullet + " Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
$d.Paragraphs($bulletBase + 2).Range.Text = $bullet + " Real-time collaboration at national scale"
$d.Paragraphs($bulletBase + 3).Range.Text = $bullet + " Revenue generation: Delivered `$4.9M additional revenue through optimization"

# Bullets 4 and 5 (Trigonometric algorithm / Discovered systematic race
# coding errors) are dropped entirely - delete the whole paragraphs
# (including their paragraph marks). Delete the later one first so the
# earlier paragraph's index doesn't shift.
$d.Paragraphs($bulletBase + 5).Range.Delete()
$d.Paragraphs($bulletBase + 4).Range.Delete()

# The remaining bullet (originally #6) becomes the new final bullet.
$d.Paragraphs($bulletBase + 4).Range.Text = $bullet + " 23% conversion rate improvement"
